# Corrections to Module Outline and Lab 2
# Applies the text corrections on the "forR" sheet's module-outline /
# lab schedule table, clears the old "Project Due Date:" label cell and
# writes the new label in column D, and restores the print orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("forR")
$ws.Activate()

# --- Topic (column D) corrections -----------------------------------------
$ws.Range("D5").Value  = "Spatial Data Types: vectors and rasters"
$ws.Range("D7").Value  = "Working with vector data - attribute tables (Summative Assessment 1)"
$ws.Range("D8").Value  = "Spatial Data Types: vectors and rasters"
$ws.Range("D9").Value  = "Spatial operations with vectors - geoprocessing"
$ws.Range("D10").Value = "Spatial operations with rasters - raster calculator and map algebra"
$ws.Range("D11").Value = "Raster terrain analysis / raster algebra"
$ws.Range("D12").Value = "Raster terrain analysis   "
$ws.Range("D13").Value = "Raster spatial operations (Summative Assessment 2)"
$ws.Range("D14").Value = "Zonal statistics / interpolation / cost analysis "
$ws.Range("D15").Value = "Combined spatial analysis"
$ws.Range("D16").Value = "Combined spatial analysis (Summative Assessment 3)"
$ws.Range("D20").Value = "Combined exercises, all previous topics (Summative Assessment 4)"
$ws.Range("D24").Value = "Remote sensing image visualization / digitization (Summative assessment 5)"
$ws.Range("D29").Value = "Satellite image classification (Summative Assessment 6 -  project proposal)"
$ws.Range("D30").Value = "Data collection and classification validation "

# --- Project due-date row (35): move the label from B to D -----------------
$ws.Range("B35").ClearContents()
$ws.Range("D35").Value = "Project Due Date"

# --- View state: scroll / selection on the forR sheet -----------------------
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D41").Select()

# --- Page setup ---------------------------------------------------------
$ws.PageSetup.Orientation = 1
